# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets
# to reflect newly generated counts (gh-pages output regeneration).

$wb = $excel.ActiveWorkbook

# Map of worksheet name -> (row -> new value) for column F updates
$updates = @{
    "展览" = @{
        5  = 13817
        7  = 249
        16 = 13863
        19 = 14932
        21 = 8225
        30 = 1033
        35 = 2
        36 = 9
        41 = 5070
    }
    "全部类型" = @{
        5  = 13817
        7  = 249
        16 = 13863
        19 = 14932
        21 = 8225
        30 = 1033
        37 = 2
        38 = 9
        43 = 5070
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowsMap = $updates[$sheetName]
    foreach ($row in $rowsMap.Keys) {
        $ws.Range("F$row").Value = $rowsMap[$row]
    }
}
